$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the "COMPLETED" format (green fill) from F8 and apply it to F9:F12
$ws.Range("F8").Copy()
$ws.Range("F9:F12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Mark rows 9-12 (S/N 8-11) as COMPLETED in the progress column
$ws.Range("F9").Value = "COMPLETED"
$ws.Range("F10").Value = "COMPLETED"
$ws.Range("F11").Value = "COMPLETED"
$ws.Range("F12").Value = "COMPLETED"

# Update the active selection to F5
$ws.Range("F5").Select()
